$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 321335
$ws.Range("D2").Value = 409522759
$ws.Range("C8").Value = 863
$ws.Range("D8").Value = 1269295
$ws.Range("C10").Value = 117147
$ws.Range("D10").Value = 171657107
$ws.Range("C12").Value = 59516
$ws.Range("D12").Value = 85899912
$ws.Range("C16").Value = 4007
$ws.Range("D16").Value = 5685992
$ws.Range("C20").Value = 6668
$ws.Range("D20").Value = 9306076
$ws.Range("C22").Value = 77629
$ws.Range("D22").Value = 96800087
$ws.Range("C27").Value = 287
$ws.Range("D27").Value = 411647
$ws.Range("C28").Value = 32496
$ws.Range("D28").Value = 47572675
$ws.Range("C30").Value = 11504
$ws.Range("D30").Value = 16548031
$ws.Range("C35").Value = 1830
$ws.Range("D35").Value = 2583295
$ws.Range("C36").Value = 97314
$ws.Range("D36").Value = 122468812
$ws.Range("C44").Value = 44412
$ws.Range("D44").Value = 65086668
$ws.Range("C46").Value = 9160
$ws.Range("D46").Value = 13145853
$ws.Range("C48").Value = 1407
$ws.Range("D48").Value = 1953103
$ws.Range("C51").Value = 2308
$ws.Range("D51").Value = 3221915
$ws.Range("C52").Value = 69166
$ws.Range("D52").Value = 86754877
$ws.Range("C59").Value = 28217
$ws.Range("D59").Value = 41380960
$ws.Range("C62").Value = 11144
$ws.Range("D62").Value = 16113338
$ws.Range("C68").Value = 1474
$ws.Range("D68").Value = 2064149
$ws.Range("C70").Value = 20537
$ws.Range("D70").Value = 26904948
$ws.Range("C74").Value = 7596
$ws.Range("D74").Value = 11123008
$ws.Range("C76").Value = 5133
$ws.Range("D76").Value = 7452515
$ws.Range("C77").Value = 492
$ws.Range("D77").Value = 696739
$ws.Range("C78").Value = 279
$ws.Range("D78").Value = 392173
$ws.Range("C79").Value = 141110
$ws.Range("D79").Value = 175974123
$ws.Range("C80").Value = 69
$ws.Range("D80").Value = 83285
$ws.Range("C83").Value = 429
$ws.Range("D83").Value = 626324
$ws.Range("C85").Value = 63596
$ws.Range("D85").Value = 93206179
$ws.Range("C88").Value = 29780
$ws.Range("D88").Value = 43080712
$ws.Range("C90").Value = 2737
$ws.Range("D90").Value = 3940857
$ws.Range("C91").Value = 2831
$ws.Range("D91").Value = 4002218
$ws.Range("C92").Value = 33168
$ws.Range("D92").Value = 44943381
$ws.Range("C96").Value = 8025
$ws.Range("D96").Value = 11800209
$ws.Range("C98").Value = 7380
$ws.Range("D98").Value = 10706154
$ws.Range("C100").Value = 533
$ws.Range("D100").Value = 757716
$ws.Range("C102").Value = 9812
$ws.Range("D102").Value = 14545354
$ws.Range("C103").Value = 5
$ws.Range("D103").Value = 7500
$ws.Range("C104").Value = 2478
$ws.Range("D104").Value = 3936378
$ws.Range("C106").Value = 3314
$ws.Range("D106").Value = 5232274
$ws.Range("C108").Value = 146
$ws.Range("D108").Value = 228820
$ws.Range("C109").Value = 187
$ws.Range("D109").Value = 275043
$ws.Range("C110").Value = 141790
$ws.Range("D110").Value = 175350331
$ws.Range("C114").Value = 950
$ws.Range("D114").Value = 1394036
$ws.Range("C116").Value = 52769
$ws.Range("D116").Value = 77349493
$ws.Range("C118").Value = 27117
$ws.Range("D118").Value = 39287777
$ws.Range("C119").Value = 1312
$ws.Range("D119").Value = 1794284
$ws.Range("C122").Value = 2265
$ws.Range("D122").Value = 3181643
$ws.Range("C124").Value = 512941
$ws.Range("D124").Value = 677422128
$ws.Range("C129").Value = 1378
$ws.Range("D129").Value = 2042486
$ws.Range("C131").Value = 208249
$ws.Range("D131").Value = 306137580
$ws.Range("C132").Value = 402
$ws.Range("D132").Value = 599750
$ws.Range("C134").Value = 183058
$ws.Range("D134").Value = 266179685
$ws.Range("C137").Value = 2852
$ws.Range("D137").Value = 4007236
$ws.Range("C139").Value = 6363
$ws.Range("D139").Value = 8986191
$ws.Range("C142").Value = 44773
$ws.Range("D142").Value = 59780360
$ws.Range("C148").Value = 14092
$ws.Range("D148").Value = 20664361
$ws.Range("C149").Value = 3761
$ws.Range("D149").Value = 5424751
$ws.Range("C152").Value = 402
$ws.Range("D152").Value = 578216
$ws.Range("C154").Value = 385
$ws.Range("D154").Value = 544163
$ws.Range("C155").Value = 17640
$ws.Range("D155").Value = 23319922
$ws.Range("C159").Value = 7189
$ws.Range("D159").Value = 10456701
$ws.Range("C161").Value = 5017
$ws.Range("D161").Value = 7220397
$ws.Range("C164").Value = 267
$ws.Range("D164").Value = 381864
$ws.Range("C166").Value = 17526
$ws.Range("D166").Value = 27658559
$ws.Range("C167").Value = 1945
$ws.Range("D167").Value = 3105649
$ws.Range("C168").Value = 260
$ws.Range("D168").Value = 412433
$ws.Range("C170").Value = 60
$ws.Range("D170").Value = 100190
$ws.Range("C172").Value = 87893
$ws.Range("D172").Value = 109911146
$ws.Range("C177").Value = 642
$ws.Range("D177").Value = 946348
$ws.Range("C179").Value = 33903
$ws.Range("D179").Value = 49720524
$ws.Range("C181").Value = 13062
$ws.Range("D181").Value = 18873514
$ws.Range("C183").Value = 1247
$ws.Range("D183").Value = 1745429
$ws.Range("C185").Value = 1652
$ws.Range("D185").Value = 2322234
$ws.Range("C186").Value = 4
$ws.Range("D186").Value = 6000
$ws.Range("C187").Value = 239039
$ws.Range("D187").Value = 297152156
$ws.Range("C193").Value = 878
$ws.Range("D193").Value = 1291345
$ws.Range("C195").Value = 86687
$ws.Range("D195").Value = 127073143
$ws.Range("C196").Value = 94
$ws.Range("D196").Value = 136627
$ws.Range("C198").Value = 33105
$ws.Range("D198").Value = 47653114
$ws.Range("C201").Value = 5114
$ws.Range("D201").Value = 7279777
$ws.Range("C204").Value = 4880
$ws.Range("D204").Value = 6757428
$ws.Range("C207").Value = 264478
$ws.Range("D207").Value = 327325610
$ws.Range("C214").Value = 617
$ws.Range("D214").Value = 898378
$ws.Range("C216").Value = 95230
$ws.Range("D216").Value = 139321251
$ws.Range("C217").Value = 89
$ws.Range("D217").Value = 132699
$ws.Range("C219").Value = 51568
$ws.Range("D219").Value = 74541518
$ws.Range("C220").Value = 34
$ws.Range("D220").Value = 48922
$ws.Range("C222").Value = 4673
$ws.Range("D222").Value = 6560805
$ws.Range("C225").Value = 5771
$ws.Range("D225").Value = 7983779
$ws.Range("C228").Value = 106574
$ws.Range("D228").Value = 133275394
$ws.Range("C233").Value = 566
$ws.Range("D233").Value = 826939
$ws.Range("C235").Value = 49559
$ws.Range("D235").Value = 72599734
$ws.Range("C237").Value = 12443
$ws.Range("D237").Value = 17892635
$ws.Range("C239").Value = 1896
$ws.Range("D239").Value = 2717882
$ws.Range("C241").Value = 2515
$ws.Range("D241").Value = 3518971
$ws.Range("C242").Value = 258296
$ws.Range("D242").Value = 326191127
$ws.Range("C248").Value = 832
$ws.Range("D248").Value = 1222063
$ws.Range("C250").Value = 95917
$ws.Range("D250").Value = 140544736
$ws.Range("C253").Value = 65239
$ws.Range("D253").Value = 94546964
$ws.Range("C255").Value = 2413
$ws.Range("D255").Value = 3403509
$ws.Range("C258").Value = 4620
$ws.Range("D258").Value = 6489371
